$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 36482.07
$ws.Range("J17").Value = 36482.07
$ws.Range("L17").Value = 109446.21
$ws.Range("N17").Value = -109782.21

$ws.Range("H129").Value = 982.71155
$ws.Range("I129").Value = 700
$ws.Range("J129").Value = 994.02
$ws.Range("K129").Value = 2100
$ws.Range("L129").Value = 2982.06
$ws.Range("M129").Value = 2900
$ws.Range("N129").Value = -12982.06

$ws.Range("H134").Value = 92852.73
$ws.Range("J134").Value = 92852.73
$ws.Range("L134").Value = 92852.73
$ws.Range("N134").Value = -102992.73

$ws.Range("H137").Value = 1363.4038
$ws.Range("I137").Value = 1187.5454
$ws.Range("J137").Value = 1668.8422
$ws.Range("K137").Value = 3562.6362
$ws.Range("L137").Value = 5006.5266
$ws.Range("M137").Value = -1012.6362
$ws.Range("N137").Value = -10106.5266

$ws.Range("H138").Value = 4308.839
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4308.839
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 12926.517
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -23206.517

$ws.Range("H140").Value = 92496.664
$ws.Range("J140").Value = 119995
$ws.Range("L140").Value = 119995
$ws.Range("N140").Value = -130355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10499.701
$ws.Range("I32").Value = 10662.237
$ws.Range("K32").Value = 10662.237
$ws.Range("M32").Value = -10375.237

$ws.Range("H45").Value = 1443.2
$ws.Range("I45").Value = 1175.909
$ws.Range("J45").Value = 2178.25
$ws.Range("K45").Value = 1175.909
$ws.Range("L45").Value = 2178.25
$ws.Range("M45").Value = -798.9090000000001
$ws.Range("N45").Value = -2932.25

$ws.Range("H55").Value = 47900
$ws.Range("J55").Value = 47900
$ws.Range("L55").Value = 47900
$ws.Range("N55").Value = -48530

$ws.Range("H61").Value = 1304
$ws.Range("I61").Value = 1081
$ws.Range("J61").Value = 3088
$ws.Range("K61").Value = 1081
$ws.Range("L61").Value = 3088
$ws.Range("M61").Value = -869
$ws.Range("N61").Value = -3512

$ws.Range("H74").Value = 1307.6666
$ws.Range("I74").Value = 1403.591
$ws.Range("J74").Value = 885.6
$ws.Range("K74").Value = 1403.591
$ws.Range("L74").Value = 885.6
$ws.Range("M74").Value = -529.5909999999999
$ws.Range("N74").Value = -2633.6

$ws.Range("H77").Value = 1307.6666
$ws.Range("I77").Value = 1403.591
$ws.Range("J77").Value = 885.6
$ws.Range("K77").Value = 7017.955
$ws.Range("L77").Value = 4428
$ws.Range("M77").Value = -2649.955
$ws.Range("N77").Value = -13164

$ws.Range("H110").Value = 1471
$ws.Range("I110").Value = 1462.8235
$ws.Range("J110").Value = 1505.75
$ws.Range("K110").Value = 1462.8235
$ws.Range("L110").Value = 1505.75
$ws.Range("M110").Value = 582.1765
$ws.Range("N110").Value = -5595.75

$ws.Range("H136").Value = 1304
$ws.Range("I136").Value = 1081
$ws.Range("J136").Value = 3088
$ws.Range("K136").Value = 3243
$ws.Range("L136").Value = 9264
$ws.Range("M136").Value = -693
$ws.Range("N136").Value = -14364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 51238.15
$ws.Range("I94").Value = 427
$ws.Range("J94").Value = 127454.875
$ws.Range("K94").Value = 427
$ws.Range("L94").Value = 127454.875
$ws.Range("M94").Value = 24
$ws.Range("N94").Value = -128356.875

$ws.Range("H105").Value = 2914.1875
$ws.Range("I105").Value = 2724.0715
$ws.Range("J105").Value = 4245
$ws.Range("K105").Value = 2724.0715
$ws.Range("L105").Value = 4245
$ws.Range("M105").Value = -977.0715
$ws.Range("N105").Value = -7739

$ws.Range("H109").Value = 21575.7
$ws.Range("J109").Value = 21575.7
$ws.Range("L109").Value = 21575.7
$ws.Range("N109").Value = -24349.7

$ws.Range("H134").Value = 3255.2068
$ws.Range("I134").Value = 5373.143
$ws.Range("J134").Value = 2581.318
$ws.Range("K134").Value = 16119.429
$ws.Range("L134").Value = 7743.954000000001
$ws.Range("M134").Value = -13584.429
$ws.Range("N134").Value = -12813.954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2242.85
$ws.Range("I31").Value = 1647.3334
$ws.Range("K31").Value = 1647.3334
$ws.Range("M31").Value = -1352.3334

$ws.Range("H34").Value = 2242.85
$ws.Range("I34").Value = 1647.3334
$ws.Range("K34").Value = 1647.3334
$ws.Range("M34").Value = -1445.3334

$ws.Range("H58").Value = 741990.3
$ws.Range("I58").Value = 1123458.2
$ws.Range("J58").Value = 1493.5883
$ws.Range("K58").Value = 1123458.2
$ws.Range("L58").Value = 1493.5883
$ws.Range("M58").Value = -1123255.2
$ws.Range("N58").Value = -1899.5883

$ws.Range("H132").Value = 256351.94
$ws.Range("I132").Value = 322931.25
$ws.Range("K132").Value = 968793.75
$ws.Range("M132").Value = -966263.75

$ws.Range("H136").Value = 741990.3
$ws.Range("I136").Value = 1123458.2
$ws.Range("J136").Value = 1493.5883
$ws.Range("K136").Value = 3370374.6
$ws.Range("L136").Value = 4480.7649
$ws.Range("M136").Value = -3367824.6
$ws.Range("N136").Value = -9580.7649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 95031.25
$ws.Range("I97").Value = 34700
$ws.Range("J97").Value = 1000000
$ws.Range("K97").Value = 34700
$ws.Range("L97").Value = 1000000
$ws.Range("M97").Value = -34204
$ws.Range("N97").Value = -1000992

$ws.Range("H107").Value = 6905.5884
$ws.Range("J107").Value = 865.1429000000001
$ws.Range("L107").Value = 865.1429000000001
$ws.Range("N107").Value = -4705.1429

$ws.Range("H122").Value = 3732
$ws.Range("I122").Value = 1900
$ws.Range("K122").Value = 5700
$ws.Range("M122").Value = -3250

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3356.2666
$ws.Range("I132").Value = 3073.3914
$ws.Range("J132").Value = 4285.7144
$ws.Range("K132").Value = 9220.174199999999
$ws.Range("L132").Value = 12857.1432
$ws.Range("M132").Value = -6690.174199999999
$ws.Range("N132").Value = -17917.1432

$ws.Range("H136").Value = 2046.46
$ws.Range("I136").Value = 1704.8445
$ws.Range("K136").Value = 5114.5335
$ws.Range("M136").Value = -2564.5335

$ws.Range("H139").Value = 56833.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 56833.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 56833.332
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -67113.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1369.7435
$ws.Range("I132").Value = 1285.8077
$ws.Range("J132").Value = 1537.6154
$ws.Range("K132").Value = 3857.4231
$ws.Range("L132").Value = 4612.8462
$ws.Range("M132").Value = -1327.4231
$ws.Range("N132").Value = -9672.8462

$ws.Range("H136").Value = 1828.3
$ws.Range("I136").Value = 1828.3
$ws.Range("K136").Value = 5484.9
$ws.Range("M136").Value = -2934.9
